$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the method under test in row 4 (A4)
$ws.Range("A4").Value = "verifyCustomerNewSeasonalLicensePurchase"

# Remove row 7 entirely (verifyCustomerReceivedSubscriptionUpgradeReceipt)
$ws.Rows.Item(7).Delete()
